$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.5254717992426
$ws.Range("D2").Value = 3.05306225987883
$ws.Range("E2").Value = 19.46157463727772
$ws.Range("F2").Value = 17.93764801243444
$ws.Range("G2").Value = 19.13130549794963
$ws.Range("H2").Value = 11.05197927298044
$ws.Range("I2").Value = 22.04651378314569
$ws.Range("L2").Value = 8.831106694377269
$ws.Range("N2").Value = 17.71690969420926
$ws.Range("O2").Value = 15.59750735714329
$ws.Range("B3").Value = 14.21647929673234
$ws.Range("D3").Value = 3.024369112635032
$ws.Range("E3").Value = 19.45830949409602
$ws.Range("F3").Value = 17.71864379238568
$ws.Range("G3").Value = 18.61633113051755
$ws.Range("H3").Value = 11.04995445959457
$ws.Range("I3").Value = 22.19105478300767
$ws.Range("L3").Value = 8.686500625716381
$ws.Range("N3").Value = 17.70234755475438
$ws.Range("O3").Value = 15.50875811661337
$ws.Range("B4").Value = 14.02480320129745
$ws.Range("D4").Value = 3.007117189528183
$ws.Range("E4").Value = 19.45859612847356
$ws.Range("F4").Value = 17.58866452936656
$ws.Range("G4").Value = 18.30018959595875
$ws.Range("H4").Value = 11.05105097273099
$ws.Range("I4").Value = 22.28438800424045
$ws.Range("L4").Value = 8.597589588567175
$ws.Range("N4").Value = 17.6954752632378
$ws.Range("O4").Value = 15.4585418422201
$ws.Range("B5").Value = 13.94630526760923
$ws.Range("D5").Value = 3.000185819757807
$ws.Range("E5").Value = 19.4592935324771
$ws.Range("F5").Value = 17.53689466272931
$ws.Range("G5").Value = 18.17161121209857
$ws.Range("H5").Value = 11.05208698005524
$ws.Range("I5").Value = 22.32357803814756
$ws.Range("L5").Value = 8.561369953207407
$ws.Range("N5").Value = 17.69319873371636
$ws.Range("O5").Value = 15.43917232788376
$ws.Range("B6").Value = 13.9332503654495
$ws.Range("D6").Value = 2.999041045819236
$ws.Range("E6").Value = 19.45944452808715
$ws.Range("F6").Value = 17.52837260852117
$ws.Range("G6").Value = 18.15028320691249
$ws.Range("H6").Value = 11.05229461024051
$ws.Range("I6").Value = 22.33015542618401
$ws.Range("L6").Value = 8.55535772488623
$ws.Range("N6").Value = 17.69285247071977
$ws.Range("O6").Value = 15.43602262265837
$ws.Range("B7").Value = 14.02374598465645
$ws.Range("D7").Value = 3.007023301032878
$ws.Range("E7").Value = 19.45860317737669
$ws.Range("F7").Value = 17.58796140458206
$ws.Range("G7").Value = 18.29845420217027
$ws.Range("H7").Value = 11.05106255833384
$ws.Range("I7").Value = 22.28491184994446
$ws.Range("L7").Value = 8.597101010637109
$ws.Range("N7").Value = 17.69544243476803
$ws.Range("O7").Value = 15.45827616513166
$ws.Range("B8").Value = 14.41939714985156
$ws.Range("D8").Value = 3.043097163155605
$ws.Range("E8").Value = 19.45997590597252
$ws.Range("F8").Value = 17.86124538923718
$ws.Range("G8").Value = 18.95387093810944
$ws.Range("H8").Value = 11.05079608512525
$ws.Range("I8").Value = 22.095401938398
$ws.Range("L8").Value = 8.781297579285912
$ws.Range("N8").Value = 17.71146155744707
$ws.Range("O8").Value = 15.56603093077027
$ws.Range("B9").Value = 15.1754529499095
$ws.Range("D9").Value = 3.116427798039722
$ws.Range("E9").Value = 19.48065123886015
$ws.Range("F9").Value = 18.42964792275113
$ws.Range("G9").Value = 20.22962326791431
$ws.Range("H9").Value = 11.0687897569815
$ws.Range("I9").Value = 21.76000274028024
$ws.Range("L9").Value = 9.139678396628149
$ws.Range("N9").Value = 17.75912097408643
$ws.Range("O9").Value = 15.81035511596174
$ws.Range("B10").Value = 15.71339214678815
$ws.Range("D10").Value = 3.171456909187281
$ws.Range("E10").Value = 19.50653120089677
$ws.Range("F10").Value = 18.86274601309827
$ws.Range("G10").Value = 21.14864087778435
$ws.Range("H10").Value = 11.09321357186028
$ws.Range("I10").Value = 21.53546473132115
$ws.Range("L10").Value = 9.398839219717207
$ws.Range("N10").Value = 17.80380385494282
$ws.Range("O10").Value = 16.00865044564308
$ws.Range("B11").Value = 15.95325421613637
$ws.Range("D11").Value = 3.196649884570983
$ws.Range("E11").Value = 19.52056450677877
$ws.Range("F11").Value = 19.06220495620209
$ws.Range("G11").Value = 21.5603929766888
$ws.Range("H11").Value = 11.10673252049513
$ws.Range("I11").Value = 21.43802369373613
$ws.Range("L11").Value = 9.515352527466357
$ws.Range("N11").Value = 17.82617634194222
$ws.Range("O11").Value = 16.10262408007487
$ws.Range("B12").Value = 16.04330602051619
$ws.Range("D12").Value = 3.206205139178124
$ws.Range("E12").Value = 19.52619814974155
$ws.Range("F12").Value = 19.13800810935972
$ws.Range("G12").Value = 21.71522817968624
$ws.Range("H12").Value = 11.11219542454822
$ws.Range("I12").Value = 21.40179824946419
$ws.Range("L12").Value = 9.559236145938621
$ws.Range("N12").Value = 17.83493750799584
$ws.Range("O12").Value = 16.13872284948739
$ws.Range("B13").Value = 16.02394762525669
$ws.Range("D13").Value = 3.204146690293128
$ws.Range("E13").Value = 19.52497072182633
$ws.Range("F13").Value = 19.12167166964035
$ws.Range("G13").Value = 21.68193264703653
$ws.Range("H13").Value = 11.11100365756046
$ws.Range("I13").Value = 21.4095701440174
$ws.Range("L13").Value = 9.549796184916683
$ws.Range("N13").Value = 17.83303785424869
$ws.Range("O13").Value = 16.13092598830039
$ws.Range("B14").Value = 15.9606788861173
$ws.Range("D14").Value = 3.197435740376173
$ws.Range("E14").Value = 19.52102161810019
$ws.Range("F14").Value = 19.0684363146704
$ws.Range("G14").Value = 21.57315410344261
$ws.Range("H14").Value = 11.1071750882376
$ws.Range("I14").Value = 21.435029925625
$ws.Range("L14").Value = 9.518967834604529
$ws.Range("N14").Value = 17.82689136196479
$ws.Range("O14").Value = 16.10558383530833
$ws.Range("B15").Value = 15.92182121885418
$ws.Range("D15").Value = 3.193326842309347
$ws.Range("E15").Value = 19.51864412488706
$ws.Range("F15").Value = 19.03586128011202
$ws.Range("G15").Value = 21.50637758605516
$ws.Range("H15").Value = 11.10487463439399
$ws.Range("I15").Value = 21.45071238478446
$ws.Range("L15").Value = 9.500052532343163
$ws.Range("N15").Value = 17.82316396688064
$ws.Range("O15").Value = 16.09012696784538
$ws.Range("B16").Value = 15.6976119270008
$ws.Range("D16").Value = 3.169813041290602
$ws.Range("E16").Value = 19.50565910887518
$ws.Range("F16").Value = 18.84975327427371
$ws.Range("G16").Value = 21.12158926099771
$ws.Range("H16").Value = 11.09237831272934
$ws.Range("I16").Value = 21.54192711356495
$ws.Range("L16").Value = 9.391193655176561
$ws.Range("N16").Value = 17.80238250618312
$ws.Range("O16").Value = 16.00258242892555
$ws.Range("B17").Value = 15.55876505435433
$ws.Range("D17").Value = 3.155423501318324
$ws.Range("E17").Value = 19.4982682915759
$ws.Range("F17").Value = 18.73615093664236
$ws.Range("G17").Value = 20.88378522709278
$ws.Range("H17").Value = 11.08532722921472
$ws.Range("I17").Value = 21.59908676842097
$ws.Range("L17").Value = 9.324030457078807
$ws.Range("N17").Value = 17.79015427211365
$ws.Range("O17").Value = 15.94982160184781
$ws.Range("B18").Value = 15.47845209204658
$ws.Range("D18").Value = 3.147162537323171
$ws.Range("E18").Value = 19.4942303859434
$ws.Range("F18").Value = 18.67104424453986
$ws.Range("G18").Value = 20.74642208179204
$ws.Range("H18").Value = 11.08149851265589
$ws.Range("I18").Value = 21.63240627287687
$ws.Range("L18").Value = 9.285272267218295
$ws.Range("N18").Value = 17.78331380890491
$ws.Range("O18").Value = 15.91983144052128
$ws.Range("B19").Value = 15.45118446082536
$ws.Range("D19").Value = 3.144368420235713
$ws.Range("E19").Value = 19.49289998532013
$ws.Range("F19").Value = 18.64904289818918
$ws.Range("G19").Value = 20.69981856928893
$ws.Range("H19").Value = 11.08024122028505
$ws.Range("I19").Value = 21.64376382004746
$ws.Range("L19").Value = 9.272128694890439
$ws.Range("N19").Value = 17.78103102531895
$ws.Range("O19").Value = 15.90973937312042
$ws.Range("B20").Value = 15.57359292027292
$ws.Range("D20").Value = 3.156953742834891
$ws.Range("E20").Value = 19.49903304136393
$ws.Range("F20").Value = 18.74822041521949
$ws.Range("G20").Value = 20.90916171254697
$ws.Range("H20").Value = 11.08605436671434
$ws.Range("I20").Value = 21.59295621990356
$ws.Range("L20").Value = 9.331193598044033
$ws.Range("N20").Value = 17.79143606055111
$ws.Range("O20").Value = 15.95540139422002
$ws.Range("B21").Value = 15.97928420471064
$ws.Range("D21").Value = 3.199406557599904
$ws.Range("E21").Value = 19.52217293814552
$ws.Range("F21").Value = 19.08406605755291
$ws.Range("G21").Value = 21.60513583024541
$ws.Range("H21").Value = 11.1082903311097
$ws.Range("I21").Value = 21.42753351825825
$ws.Range("L21").Value = 9.528029614033036
$ws.Range("N21").Value = 17.82868892814058
$ws.Range("O21").Value = 16.11301375246137
$ws.Range("B22").Value = 16.23985394097009
$ws.Range("D22").Value = 3.227237095807417
$ws.Range("E22").Value = 19.53915603026623
$ws.Range("F22").Value = 19.30511262605257
$ws.Range("G22").Value = 22.05359070760298
$ws.Range("H22").Value = 11.12482389559121
$ws.Range("I22").Value = 21.32334379174152
$ws.Range("L22").Value = 9.655272560116526
$ws.Range("N22").Value = 17.85471899650638
$ws.Range("O22").Value = 16.2189985928276
$ws.Range("B23").Value = 16.10122684179249
$ws.Range("D23").Value = 3.212378171801747
$ws.Range("E23").Value = 19.52992351747973
$ws.Range("F23").Value = 19.18701956742913
$ws.Range("G23").Value = 21.81488298362617
$ws.Range("H23").Value = 11.11581752255645
$ws.Range("I23").Value = 21.37859373344852
$ws.Range("L23").Value = 9.587501174613951
$ws.Range("N23").Value = 17.84067395646092
$ws.Range("O23").Value = 16.16216996738868
$ws.Range("B24").Value = 15.56689075597368
$ws.Range("D24").Value = 3.156261883387776
$ws.Range("E24").Value = 19.49868663980863
$ws.Range("F24").Value = 18.74276315998602
$ws.Range("G24").Value = 20.89769100321339
$ws.Range("H24").Value = 11.08572492663087
$ws.Range("I24").Value = 21.59572641510443
$ws.Range("L24").Value = 9.327955591206946
$ws.Range("N24").Value = 17.79085597268724
$ws.Range("O24").Value = 15.95287770027522
$ws.Range("B25").Value = 14.97362378326025
$ws.Range("D25").Value = 3.096356339245924
$ws.Range("E25").Value = 19.47316035379485
$ws.Range("F25").Value = 18.2728429432021
$ws.Range("G25").Value = 19.8868573731673
$ws.Range("H25").Value = 11.06194702992836
$ws.Range("I25").Value = 21.8468799857447
$ws.Range("L25").Value = 9.043284010725545
$ws.Range("N25").Value = 17.74451223896402
$ws.Range("O25").Value = 15.81035511596174
